# Add two new "Model" result blocks (Model 2.6 @ rows 181-187, Model 2.7 @
# rows 189-195) to the results table on sheet1, mirroring the existing
# "Model 2.x" blocks already present (e.g. rows 165-171 / 173-179).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Bring over cell formatting (styles / column-level formats) for the
#    two new 7-row blocks by copying the formats of the most recent
#    existing blocks (rows 165-171 and 173-179) onto the new rows.
# ---------------------------------------------------------------------
$srcRows1 = 165
$dstRows1 = 181
$srcRows2 = 173
$dstRows2 = 189

for ($i = 0; $i -lt 7; $i++) {
    $srcRow = $srcRows1 + $i
    $dstRow = $dstRows1 + $i
    $ws.Range("A" + $srcRow + ":T" + $srcRow).Copy()
    $ws.Range("A" + $dstRow + ":T" + $dstRow).PasteSpecial(-4122)
}

for ($i = 0; $i -lt 7; $i++) {
    $srcRow = $srcRows2 + $i
    $dstRow = $dstRows2 + $i
    $ws.Range("A" + $srcRow + ":T" + $srcRow).Copy()
    $ws.Range("A" + $dstRow + ":T" + $dstRow).PasteSpecial(-4122)
}

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 2) Block 1 -- "Model 2.6" (rows 181-187)
# ---------------------------------------------------------------------
$ws.Range("A181").Value = "Model 2.6"
$ws.Range("B181").Value = "(64, 64, 1)"
$ws.Range("C181").Value = 32
$ws.Range("I181").Value = 0.0001
$ws.Range("T181").Formula = "=L181-J181"

$ws.Range("A182").Value = "Augmentation"
$ws.Range("O182").Value = 0.5
$ws.Range("P182").Value = 0.2
$ws.Range("Q182").Value = "(0,3 0,3)"
$ws.Range("R182").Value = "hor_ver"

$ws.Range("A183").Value = "Conv 1"
$ws.Range("D183").Value = 32
$ws.Range("E183").Value = "(3,3)"
$ws.Range("F183").Value = "(2,2)"
$ws.Range("G183").Value = "no"
$ws.Range("H183").Value = 0.2
$ws.Range("T183").Formula = "=L183-J183"

$ws.Range("A184").Value = "Conv 2"
$ws.Range("D184").Value = 64
$ws.Range("E184").Value = "(3,3)"
$ws.Range("F184").Value = "(2,2)"
$ws.Range("G184").Value = "no"
$ws.Range("H184").Value = 0.3
$ws.Range("T184").Formula = "=L184-J184"

$ws.Range("A185").Value = "Conv 3"
$ws.Range("D185").Value = 64
$ws.Range("E185").Value = "(5,5)"
$ws.Range("F185").Value = "(2,2)"
$ws.Range("G185").Value = "no"
$ws.Range("H185").Value = 0.5
$ws.Range("T185").Formula = "=L185-J185"

$ws.Range("A186").Value = "Dense"
$ws.Range("D186").Value = 128
$ws.Range("H186").Value = 0.5
$ws.Range("T186").Formula = "=L186-J186"

$ws.Range("A187").Value = "Output"
$ws.Range("J187").Value = 0.5753
$ws.Range("K187").Value = 0.7021
$ws.Range("L187").Value = 2.5042
$ws.Range("M187").Value = 0.4611
$ws.Range("N187").Value = 35
$ws.Range("T187").Formula = "=L187-J187"

# ---------------------------------------------------------------------
# 3) Block 2 -- "Model 2.7" (rows 189-195)
# ---------------------------------------------------------------------
$ws.Range("A189").Value = "Model 2.7"
$ws.Range("B189").Value = "(64, 64, 1)"
$ws.Range("C189").Value = 16
$ws.Range("I189").Value = 0.0001
$ws.Range("T189").Formula = "=L189-J189"

$ws.Range("A190").Value = "Augmentation"
$ws.Range("O190").Value = 0.2
$ws.Range("P190").Value = 0.2
$ws.Range("Q190").Value = "(0,1 0,1)"
$ws.Range("R190").Value = "hor"

$ws.Range("A191").Value = "Conv 1"
$ws.Range("D191").Value = 32
$ws.Range("E191").Value = "(3,3)"
$ws.Range("F191").Value = "(2,2)"
$ws.Range("G191").Value = "no"
$ws.Range("H191").Value = 0.2
$ws.Range("T191").Formula = "=L191-J191"

$ws.Range("A192").Value = "Conv 2"
$ws.Range("D192").Value = 64
$ws.Range("E192").Value = "(3,3)"
$ws.Range("F192").Value = "(2,2)"
$ws.Range("G192").Value = "no"
$ws.Range("H192").Value = 0.3
$ws.Range("T192").Formula = "=L192-J192"

$ws.Range("A193").Value = "Conv 3"
$ws.Range("D193").Value = 64
$ws.Range("E193").Value = "(5,5)"
$ws.Range("F193").Value = "(2,2)"
$ws.Range("G193").Value = "no"
$ws.Range("H193").Value = 0.5
$ws.Range("T193").Formula = "=L193-J193"

$ws.Range("A194").Value = "Dense"
$ws.Range("D194").Value = 128
$ws.Range("H194").Value = 0.5
$ws.Range("T194").Formula = "=L194-J194"

$ws.Range("A195").Value = "Output"
$ws.Range("J195").Value = 0.409
$ws.Range("K195").Value = 0.8151
$ws.Range("L195").Value = 0.7457
$ws.Range("M195").Value = 0.7044
$ws.Range("N195").Value = 40
$ws.Range("T195").Formula = "=L195-J195"

# ---------------------------------------------------------------------
# 4) View state: scroll/select so the active window ends up looking at
#    the newly added data, matching the author's final selection.
# ---------------------------------------------------------------------
$ws.Range("A174").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("T195").Select()
